$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 62: 12-Dic-23, a comprobar, 500
$ws.Range("B62").Value = 45272
$ws.Range("C62").Value = "a comprobar"
$ws.Range("D62").Value = 500

# Row 63: 01-Dic-23, 4 botellones, -212
$ws.Range("B63").Value = 45261
$ws.Range("C63").Value = "4 botellones"
$ws.Range("D63").Value = -212

# Row 64: 05-Dic-23, 2 Botellones, -106
$ws.Range("B64").Value = 45265
$ws.Range("C64").Value = "2 Botellones"
$ws.Range("D64").Value = -106

# Row 65: 08-Dic-23, 3 botellones, -159
$ws.Range("B65").Value = 45268
$ws.Range("C65").Value = "3 botellones"
$ws.Range("D65").Value = -159

# Update sheet view to reflect new scroll/selection position
$ws.Application.ActiveWindow.ScrollRow = 54
$ws.Range("D63").Select()
